# Sub_regional mapping.xlsx edit script
# 1. Add a new "Inner-Outer" worksheet (alphabetical borough -> Inner/Outer lookup table)
#    positioned right after the existing "Mapping" sheet.
# 2. Add a new column C ("Inner/outer") to the "Mapping" sheet that looks up each
#    borough's Inner/Outer classification via INDEX/MATCH against the new sheet.

$wb = $excel.ActiveWorkbook
$mapping = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Create the "Normal 2" cell style (Arial) used by the new lookup sheet, mirroring
# the style that Excel creates when data is pasted in from another workbook.
# ---------------------------------------------------------------------------
$normal2 = $wb.Styles.Add("Normal 2")
$normal2.Font.Name = "Calibri"

# ---------------------------------------------------------------------------
# 1) New "Inner-Outer" worksheet, inserted directly after "Mapping"
# ---------------------------------------------------------------------------
$innerOuter = $wb.Worksheets.Add($null, $mapping)
$innerOuter.Name = "Inner-Outer"

$innerOuterData = @(
    @("Barking and Dagenham", "Outer"),
    @("Barnet", "Outer"),
    @("Bexley", "Outer"),
    @("Brent", "Outer"),
    @("Bromley", "Outer"),
    @("Camden", "Inner"),
    @("City of London", "Inner"),
    @("Croydon", "Outer"),
    @("Ealing", "Outer"),
    @("Enfield", "Outer"),
    @("Greenwich", "Inner"),
    @("Hackney", "Inner"),
    @("Hammersmith and Fulham", "Inner"),
    @("Haringey", "Outer"),
    @("Harrow", "Outer"),
    @("Havering", "Outer"),
    @("Hillingdon", "Outer"),
    @("Hounslow", "Outer"),
    @("Islington", "Inner"),
    @("Kensington and Chelsea", "Inner"),
    @("Kingston upon Thames", "Outer"),
    @("Lambeth", "Inner"),
    @("Lewisham", "Inner"),
    @("Merton", "Outer"),
    @("Newham", "Inner"),
    @("Redbridge", "Outer"),
    @("Richmond upon Thames", "Outer"),
    @("Southwark", "Inner"),
    @("Sutton", "Outer"),
    @("Tower Hamlets", "Inner"),
    @("Waltham Forest", "Outer"),
    @("Wandsworth", "Inner"),
    @("Westminster", "Inner")
)

for ($i = 0; $i -lt $innerOuterData.Count; $i++) {
    $r = $i + 1
    $innerOuter.Cells.Item($r, 1).Value = $innerOuterData[$i][0]
    $innerOuter.Cells.Item($r, 2).Value = $innerOuterData[$i][1]
}

$usedRange = $innerOuter.Range("A1:B" + $innerOuterData.Count)
$usedRange.Style = "Normal 2"
$usedRange.Font.Name = "Arial"

$innerOuter.Columns.Item(1).ColumnWidth = 22.1
$innerOuter.PageSetup.PaperSize = 9
$innerOuter.PageSetup.Orientation = 1
$innerOuter.Range("A17").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2) "Mapping" sheet: new column C with the Inner/outer lookup formula
# ---------------------------------------------------------------------------
$mapping.Range("C1").Value = "Inner/outer"
$mapping.Range("C1").Copy() | Out-Null
$mapping.Range("C1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats (copy header style from B1 would also work, reuse own)
$excel.CutCopyMode = 0

# Match the header's existing style (same as B1)
$mapping.Range("B1").Copy() | Out-Null
$mapping.Range("C1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$mapping.Range("C1").Value = "Inner/outer"

for ($r = 2; $r -le 34; $r++) {
    $mapping.Range("C$r").Formula = "=INDEX('Inner-Outer'!`$A:`$B,MATCH(`$A$r,'Inner-Outer'!A:A,0),2)"
}

$mapping.Range("B2:B34").Copy() | Out-Null
$mapping.Range("C2:C34").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$mapping.Range("C11").Select() | Out-Null

Write-Output "edit complete"
